# Update "想去人数" (number of people interested) counters.
# These values were regenerated (incremented by 1) for the gh-pages
# data refresh, both on the "展览" sheet and the mirrored "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAll        = $wb.Worksheets.Item("全部类型")

# Sheet "展览": rows 6, 9, 10, 11, 13 -> column F
$wsExhibition.Range("F6").Value  = 9881
$wsExhibition.Range("F9").Value  = 1239
$wsExhibition.Range("F10").Value = 3944
$wsExhibition.Range("F11").Value = 176
$wsExhibition.Range("F13").Value = 49

# Sheet "全部类型": rows 7, 10, 11, 12, 14 -> column F (same events, shifted by one row)
$wsAll.Range("F7").Value  = 9881
$wsAll.Range("F10").Value = 1239
$wsAll.Range("F11").Value = 3944
$wsAll.Range("F12").Value = 176
$wsAll.Range("F14").Value = 49
